$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-9 (columns A, B, C)
$data = @(
    @("<delete>", "<delete>", 52),
    @("<escape>", "<eight>", 51),
    @("<she>", "<senten>", 49),
    @("<out>", "<of>", 51),
    @("<its>", "<are>", 55),
    @("<whiskey>", "<which>", 43),
    @("<will>", "<will>", 39),
    @("<yes>", "<it>", 13)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Delete rows 10-16 which are no longer present (shrinks dimension to A1:C9)
$ws.Range("A10:C16").EntireRow.Delete()
